$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 777, shifting existing data
# (old rows 777-818) down to rows 779-820.
$ws.Rows("777:778").Insert()

# New row 777: Palta / Hass / Especial record dated 44610, Provincia de Quillota
$ws.Cells.Item(777, 1).Value2 = 5
$ws.Cells.Item(777, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(777, 3).Value2 = "Maule"
$ws.Cells.Item(777, 4).Value2 = 44610
$ws.Cells.Item(777, 5).Value2 = 7
$ws.Cells.Item(777, 6).Value2 = "Fruta"
$ws.Cells.Item(777, 7).Value2 = 100106
$ws.Cells.Item(777, 8).Value2 = "Oleaginosos"
$ws.Cells.Item(777, 9).Value2 = 100106002
$ws.Cells.Item(777, 10).Value2 = "Palta"
$ws.Cells.Item(777, 11).Value2 = "Hass"
$ws.Cells.Item(777, 12).Value2 = "Especial"
$ws.Cells.Item(777, 13).Value2 = 150
$ws.Cells.Item(777, 14).Value2 = 2700
$ws.Cells.Item(777, 15).Value2 = 2700
$ws.Cells.Item(777, 16).Value2 = 2700
$ws.Cells.Item(777, 17).Value2 = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(777, 18).Value2 = "Provincia de Quillota"
$ws.Cells.Item(777, 19).Value2 = 2700
$ws.Cells.Item(777, 20).Value2 = 1

# New row 778: Palta / Hass / Primera record dated 44610, Provincia de Quillota
$ws.Cells.Item(778, 1).Value2 = 5
$ws.Cells.Item(778, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(778, 3).Value2 = "Maule"
$ws.Cells.Item(778, 4).Value2 = 44610
$ws.Cells.Item(778, 5).Value2 = 7
$ws.Cells.Item(778, 6).Value2 = "Fruta"
$ws.Cells.Item(778, 7).Value2 = 100106
$ws.Cells.Item(778, 8).Value2 = "Oleaginosos"
$ws.Cells.Item(778, 9).Value2 = 100106002
$ws.Cells.Item(778, 10).Value2 = "Palta"
$ws.Cells.Item(778, 11).Value2 = "Hass"
$ws.Cells.Item(778, 12).Value2 = "Primera"
$ws.Cells.Item(778, 13).Value2 = 100
$ws.Cells.Item(778, 14).Value2 = 2500
$ws.Cells.Item(778, 15).Value2 = 2500
$ws.Cells.Item(778, 16).Value2 = 2500
$ws.Cells.Item(778, 17).Value2 = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(778, 18).Value2 = "Provincia de Quillota"
$ws.Cells.Item(778, 19).Value2 = 2500
$ws.Cells.Item(778, 20).Value2 = 1
